$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 2.27127949656926
$ws.Range("C2").Value = 2.27127949656926
$ws.Range("D2").Value = 2.05284250815268
$ws.Range("E2").Value = 0.0138437179182969
$ws.Range("F2").Value = 0.1567

# Row 3
$ws.Range("B3").Value = 2.22439408190873
$ws.Range("C3").Value = 2.22439408190873
$ws.Range("D3").Value = 2.01046622977176
$ws.Range("E3").Value = 0.0135579457550633
$ws.Range("F3").Value = 0.1611

# Row 4
$ws.Range("B4").Value = 0.24742240367529
$ws.Range("C4").Value = 0.24742240367529
$ws.Range("D4").Value = 0.223626915358129
$ws.Range("E4").Value = 0.00150806889610965
$ws.Range("F4").Value = 0.6524

# Row 5
$ws.Range("B5").Value = 159.322620321367
$ws.Range("C5").Value = 1.10640708556505
$ws.Range("E5").Value = 0.97109026743053

# Row 6
$ws.Range("B6").Value = 164.065716303521
